# Updated symbol list on Fri Jan 13 07:56:07 UTC 2023 with GitHub Actions
# Refresh of Price (D) and Volume(1h) (E) columns for the crypto ticker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", "287.89")
    ,@("E2", "0.96%")
    ,@("D3", "29.62")
    ,@("E3", "3.95%")
    ,@("E4", "0.62%")
    ,@("D5", "0.06687")
    ,@("E5", "3.21%")
    ,@("E6", "1.59%")
    ,@("E7", "1.22%")
    ,@("D8", "1.364")
    ,@("E8", "1.52%")
    ,@("D9", "0.9192")
    ,@("E9", "0.73%")
    ,@("D10", "0.1597")
    ,@("E10", "3.26%")
    ,@("D11", "0.06812")
    ,@("E11", "5.21%")
    ,@("D12", "0.07595")
    ,@("E12", "-0.89%")
    ,@("D13", "0.02937")
    ,@("E13", "-1.58%")
    ,@("D14", "0.08975")
    ,@("E14", "0.32%")
    ,@("D15", "0.001575")
    ,@("E15", "-1.11%")
    ,@("D16", "0.04500")
    ,@("E16", "0.85%")
    ,@("D17", "0.0006441")
    ,@("E17", "-1.61%")
    ,@("D18", "0.006294")
    ,@("E18", "3.13%")
    ,@("D19", "3.449")
    ,@("E19", "-0.28%")
    ,@("D20", "2.229")
    ,@("E20", "-0.57%")
    ,@("E21", "1.12%")
    ,@("E22", "-2.44%")
    ,@("D23", "4.079")
    ,@("E23", "2.72%")
    ,@("D24", "0.1582")
    ,@("E24", "1.72%")
    ,@("D25", "0.001187")
    ,@("E25", "0.68%")
    ,@("D26", "0.004113")
    ,@("E26", "-4.83%")
    ,@("E27", "1.52%")
    ,@("E40", "2.75%")
    ,@("D41", "0.006702")
    ,@("E41", "-1.13%")
    ,@("D42", "0.1240")
    ,@("E42", "0.74%")
    ,@("D43", "0.002228")
    ,@("E43", "6.00%")
    ,@("D44", "0.01337")
    ,@("E44", "12.30%")
    ,@("D45", "0.00005674")
    ,@("E45", "5.10%")
    ,@("D46", "1.974")
    ,@("E46", "2.14%")
    ,@("E47", "-29.41%")
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $newValue = $pair[1]
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "287.89")
    # and percentages (e.g. "0.96%") are kept as literal text,
    # matching the inline-string cells already used in this sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.NumberFormat = "General"
}
